$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the unused "Middle Name" column (column D) entirely, shifting
# everything to its right one column to the left.
$ws.Range("D1").EntireColumn.Delete()

# After deleting the column, Excel leaves the resulting column D selected.
$ws.Columns("D").Select() | Out-Null
